# regen save_data to use K instead of Strike#, regen std/mean, calc and write s_vals
# This updates the "K" column (column G) values for each row, which were
# recomputed upstream (data regen) and re-saved into the workbook.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 3
    3  = 3
    4  = 0
    5  = 0
    6  = 5
    7  = 3
    8  = 3
    9  = 2
    10 = 2
    11 = 2
    13 = 3
    14 = 3
    15 = 3
    16 = 0
    17 = 1
    18 = 1
    19 = 1
    20 = 2
    21 = 2
    22 = 2
    23 = 2
    24 = 0
    25 = 4
    26 = 2
    27 = 4
    28 = 1
    29 = 3
    30 = 1
    31 = 2
    32 = 2
    33 = 0
    34 = 5
    35 = 5
    37 = 4
    38 = 3
}

foreach ($row in $kValues.Keys) {
    $ws.Cells.Item($row, 7).Value = $kValues[$row]
}
